$wb = $excel.ActiveWorkbook

function Add-ScrimRow($ws, $formatSourceRow, $targetRow, $values, $timestamp) {
    # Copy formatting (fill/border/font) from an existing row that already
    # has the correct style for each column (including the "Equipo 1"/"Equipo 2"
    # bold fill in column G), then overwrite the cell values.
    $ws.Range("A$formatSourceRow`:N$formatSourceRow").Copy() | Out-Null
    $ws.Range("A$targetRow`:N$targetRow").PasteSpecial(-4122) | Out-Null

    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($targetRow, $i + 1).Value = $values[$i]
    }
    $ws.Cells.Item($targetRow, 14).Value = $timestamp
}

# ---------------------------------------------------------------------------
# "Triple Dribble" (sheet 1): add rows 48 and 49
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Triple Dribble")

Add-ScrimRow $ws1 47 48 @("BONNIE","BULL","CORDELIUS","LOU","HANK","LUMI","Equipo 1","HMB|Lukii","HMB|Symantec","SK|Joker","FUT|Ferissa","FUT|DeMaster","FUT|ZеyroX🕊️") "20250724T192428.000Z"
Add-ScrimRow $ws1 47 49 @("BONNIE","BULL","CORDELIUS","LOU","HANK","LUMI","Equipo 1","HMB|Lukii","HMB|Symantec","SK|Joker","FUT|Ferissa","FUT|DeMaster","FUT|ZеyroX🕊️") "20250724T192118.000Z"

# ---------------------------------------------------------------------------
# "Sneaky Fields" (sheet 10): add rows 8, 9 and 10
# ---------------------------------------------------------------------------
$ws10 = $wb.Worksheets.Item("Sneaky Fields")

Add-ScrimRow $ws10 4 8 @("LUMI","HANK","TARA","BEA","STU","EMZ","Equipo 2","FUT|GeRo","FUT|Nowy297","FUT|MeOw","Enraged 💔","SUP|Filippo神","SUP|Tomzy") "20250724T191738.000Z"
Add-ScrimRow $ws10 4 9 @("LUMI","HANK","TARA","BEA","STU","EMZ","Equipo 2","FUT|GeRo","FUT|Nowy297","FUT|MeOw","Enraged 💔","SUP|Filippo神","SUP|Tomzy") "20250724T191448.000Z"
Add-ScrimRow $ws10 6 10 @("LUMI","HANK","TARA","BEA","STU","EMZ","Equipo 1","FUT|GeRo","FUT|Nowy297","FUT|MeOw","Enraged 💔","SUP|Filippo神","SUP|Tomzy") "20250724T191208.000Z"

# ---------------------------------------------------------------------------
# "Dueling Beetles" (sheet 13): add rows 19, 20 and 21
# ---------------------------------------------------------------------------
$ws13 = $wb.Worksheets.Item("Dueling Beetles")

Add-ScrimRow $ws13 6 19 @("CORDELIUS","GRIFF","MEG","WILLOW","DRACO","KIT","Equipo 1","SK|Joker","HMB|Symantec","HMB|Lukii","FUT|ZеyroX🕊️","FUT|DeMaster","FUT|Ferissa") "20250724T191615.000Z"
Add-ScrimRow $ws13 6 20 @("CORDELIUS","GRIFF","MEG","WILLOW","DRACO","KIT","Equipo 1","SK|Joker","HMB|Symantec","HMB|Lukii","FUT|ZеyroX🕊️","FUT|DeMaster","FUT|Ferissa") "20250724T191426.000Z"
Add-ScrimRow $ws13 4 21 @("CORDELIUS","GRIFF","MEG","WILLOW","DRACO","KIT","Equipo 2","SK|Joker","HMB|Symantec","HMB|Lukii","FUT|ZеyroX🕊️","FUT|DeMaster","FUT|Ferissa") "20250724T191247.000Z"

# ---------------------------------------------------------------------------
# "Goldarm Gulch" (sheet 15): add rows 12, 13 and 14
# ---------------------------------------------------------------------------
$ws15 = $wb.Worksheets.Item("Goldarm Gulch")

Add-ScrimRow $ws15 4 12 @("CHARLIE","BUSTER","ANGELO","GUS","BONNIE","BROCK","Equipo 2","FUT|GeRo","FUT|Nowy297","FUT|MeOw","Enraged 💔","SUP|Filippo神","SUP|Tomzy") "20250724T190533.000Z"
Add-ScrimRow $ws15 5 13 @("CHARLIE","BUSTER","ANGELO","GUS","BONNIE","BROCK","Equipo 1","FUT|GeRo","FUT|Nowy297","FUT|MeOw","Enraged 💔","SUP|Filippo神","SUP|Tomzy") "20250724T190331.000Z"
Add-ScrimRow $ws15 4 14 @("CHARLIE","BUSTER","ANGELO","GUS","BONNIE","BROCK","Equipo 2","FUT|GeRo","FUT|Nowy297","FUT|MeOw","Enraged 💔","SUP|Filippo神","SUP|Tomzy") "20250724T190100.000Z"

# ---------------------------------------------------------------------------
# "Double Swoosh" (sheet 18): add rows 8 and 9
# ---------------------------------------------------------------------------
$ws18 = $wb.Worksheets.Item("Double Swoosh")

Add-ScrimRow $ws18 4 8 @("LILY","JANET","AMBER","CHARLIE","SANDY","STU","Equipo 1","SK|Joker","HMB|Lukii","HMB|Symantec","FUT|ZеyroX🕊️","FUT|Ferissa","FUT|DeMaster") "20250724T190618.000Z"
Add-ScrimRow $ws18 4 9 @("LILY","JANET","AMBER","CHARLIE","SANDY","STU","Equipo 1","SK|Joker","HMB|Lukii","HMB|Symantec","FUT|ZеyroX🕊️","FUT|Ferissa","FUT|DeMaster") "20250724T190323.000Z"
